# Updated data and cart page
#
# 1. "Sheet1" drives the generated user/email list via CONCATENATE formulas
#    that reference the absolute cell $I$23. Bumping that counter from 8 to 9
#    regenerates every "<Name>8" / "<Name>8@gmail.com" value to "<Name>9" /
#    "<Name>9@gmail.com" for rows 23-42.
$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("I23").Value = 9

# 2. The "order" sheet (rows 2-21) has its own literal (non-formula) copies of
#    those same "<Name>8" / "<Name>8@gmail.com" strings in columns R, S and T.
#    Refresh them the same way, 8 -> 9.
$order = $wb.Worksheets.Item("order")

for ($r = 2; $r -le 21; $r++) {
    $rCell = $order.Cells.Item($r, 18)  # column R
    $sCell = $order.Cells.Item($r, 19)  # column S
    $tCell = $order.Cells.Item($r, 20)  # column T

    $rVal = $rCell.Value()
    $sVal = $sCell.Value()
    $tVal = $tCell.Value()

    $rCell.Value = $rVal.Replace("8", "9")
    $sCell.Value = $sVal.Replace("8", "9")
    $tCell.Value = $tVal.Replace("8", "9")
}

# 3. Fix the first-name typo in column V row 13 ("Ty" -> "Tyth") and highlight
#    the corrected/updated first-name cells (V8, V10, V13, V19) with a yellow
#    fill, matching the other centered/bordered cells in that column.
$order.Range("V13").Value = "Tyth"

$order.Range("V8").Interior.Color = 65535
$order.Range("V10").Interior.Color = 65535
$order.Range("V13").Interior.Color = 65535
$order.Range("V19").Interior.Color = 65535

# 4. Update the saved view state on the "order" sheet: the selection moves to
#    V19 (the last cell touched above).
$order.Activate()
$order.Range("V19").Select()
